$d = $word.ActiveDocument

# Update the date heading in the first paragraph
$dateRange = $d.Paragraphs.Item(1).Range
$dateRange.Find.Execute("2024-08-14 Wednesday", $true, $false, $false, $false, $false, $true, 1, $false, "2024-08-15 Thursday", 2) | Out-Null

# New values for the 20x5 practice table, in row-major (top-to-bottom, left-to-right) order
$newValues = @(
    "66-25=", "11+39=", "46+21=", "58+12=", "29+10=",
    "62-2=", "41+38=", "94-44=", "60-24=", "12-7=",
    "99-19=", "32+65=", "75-1=", "26+71=", "12+41=",
    "62-61=", "38-8=", "42-6=", "86-1=", "59+25=",
    "53-30=", "48-10=", "35-30=", "60-18=", "44-18=",
    "17+30=", "1+40=", "83-1=", "67-18=", "22+21=",
    "9+27=", "71-65=", "92-82=", "29+36=", "94-13=",
    "88-65=", "25+32=", "60+31=", "2+71=", "61-60=",
    "24+42=", "16+69=", "75-51=", "83-56=", "21+48=",
    "9+24=", "71-57=", "42+37=", "62-29=", "46-46=",
    "43-3=", "46+50=", "71+4=", "53-36=", "50-40=",
    "12+50=", "84+7=", "40-23=", "46-3=", "89-63=",
    "21+50=", "81-14=", "47-40=", "75-11=", "56+28=",
    "46-18=", "49+50=", "35+17=", "93-67=", "27+6=",
    "98-30=", "61-46=", "26+40=", "40+17=", "29-24=",
    "47+39=", "90-63=", "11+54=", "37+60=", "8+79=",
    "81-68=", "13+49=", "8+20=", "2+38=", "14-13=",
    "33-15=", "79+7=", "82-82=", "52+26=", "70+22=",
    "69-3=", "60-6=", "69-65=", "30+22=", "55-16=",
    "27+44=", "59-54=", "40+41=", "57-26=", "3+42="
)

$t = $d.Tables.Item(1)
$rowCount = $t.Rows.Count
$colCount = $t.Columns.Count
$idx = 0
for ($r = 1; $r -le $rowCount; $r++) {
    for ($c = 1; $c -le $colCount; $c++) {
        $cell = $t.Cell($r, $c)
        $cell.Range.Text = $newValues[$idx]
        $idx++
    }
}

Write-Output ("Updated cells: " + $idx)
